$wb = $excel.ActiveWorkbook

# --- Sheet "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")

# Row 2 (Rapid Decline)
$ws1.Range("C2").Value = 3
$ws1.Range("E2").Value = 50

# Row 3 (Decline)
$ws1.Range("C3").Value = 1
$ws1.Range("E3").Value = 16.7

# Row 4 (Stable)
$ws1.Range("C4").Value = 2

# Row 5 (Increase)
$ws1.Range("B5").Value = 1
$ws1.Range("D5").Value = 50

# Row 6 (Rapid Increase)
$ws1.Range("D6").Value = 50

# Row 7 (Trend Inconclusive)
$ws1.Range("B7").Value = 29
$ws1.Range("C7").Value = 57

# --- Sheet "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")

# Row 3 (Long-term Analysis)
$ws4.Range("C3").Value = 2

# Row 4 (Current Analysis)
$ws4.Range("C4").Value = 6
